$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of the existing header cell (AC1) onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the season record values for each data row (rows 2-48)
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 86  # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 76  # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF -> Ties
}
